# Cloud ID complete through Orbits 2-44 even, 15, and 45 (Dallin sheet).
# Adds orbit-interval data rows 36-92 (orbits 22, 24, 26, 28, 30, 32, 34, 36,
# 38, 40, 42, 44 and 45) to the "Dallin" worksheet, mirroring the rows
# already present for orbits 2-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dallin")

# ---------------------------------------------------------------
# Step 1: write the new "File Name" text values first, in the same
# order they were originally typed, so the shared-string table gets
# rebuilt with entries in the exact original sequence (note orbit 34
# was entered before orbit 32).
# ---------------------------------------------------------------
$ws.Cells.Item(36, 2).Value = "awe_l1r_q20_2023326T0108_00022_v01.nc"
$ws.Cells.Item(39, 2).Value = "awe_l1r_q20_2023326T0108_00024_v01.nc"
$ws.Cells.Item(44, 2).Value = "awe_l1r_q20_2023326T0108_00026_v01.nc"
$ws.Cells.Item(49, 2).Value = "awe_l1r_q20_2023326T0108_00028_v01.nc"
$ws.Cells.Item(54, 2).Value = "awe_l1r_q20_2023326T0108_00030_v01.nc"
$ws.Cells.Item(60, 2).Value = "awe_l1r_q20_2023326T0108_00034_v01.nc"
$ws.Cells.Item(59, 2).Value = "awe_l1r_q20_2023326T0108_00032_v01.nc"
$ws.Cells.Item(69, 2).Value = "awe_l1r_q20_2023326T0108_00036_v01.nc"
$ws.Cells.Item(72, 2).Value = "awe_l1r_q20_2023326T0108_00038_v01.nc"
$ws.Cells.Item(73, 2).Value = "awe_l1r_q20_2023326T0108_00040_v01.nc"
$ws.Cells.Item(80, 2).Value = "awe_l1r_q20_2023326T0108_00042_v01.nc"
$ws.Cells.Item(86, 2).Value = "awe_l1r_q20_2023326T0108_00044_v01.nc"
$ws.Cells.Item(90, 2).Value = "awe_l1r_q20_2023326T0108_00045_v01.nc"

# ---------------------------------------------------------------
# Step 2: fill in the remaining numeric cells for every new row.
# ---------------------------------------------------------------
$ws.Cells.Item(36, 1).Value = 22
$ws.Cells.Item(36, 3).Value = 1432
$ws.Cells.Item(36, 4).Value = 558
$ws.Cells.Item(36, 5).Value = 756
$ws.Cells.Item(36, 6).Value = 596
$ws.Cells.Item(36, 7).Value = 646

$ws.Cells.Item(37, 6).Value = 695
$ws.Cells.Item(37, 7).Value = 713

$ws.Cells.Item(38, 4).Value = 992
$ws.Cells.Item(38, 5).Value = 1112
$ws.Cells.Item(38, 6).Value = 1030
$ws.Cells.Item(38, 7).Value = 1082

$ws.Cells.Item(39, 1).Value = 24
$ws.Cells.Item(39, 3).Value = 1424
$ws.Cells.Item(39, 4).Value = 948
$ws.Cells.Item(39, 5).Value = 1390
$ws.Cells.Item(39, 6).Value = 978
$ws.Cells.Item(39, 7).Value = 1035

$ws.Cells.Item(40, 6).Value = 1066
$ws.Cells.Item(40, 7).Value = 1136

$ws.Cells.Item(41, 6).Value = 1158
$ws.Cells.Item(41, 7).Value = 1158

$ws.Cells.Item(42, 6).Value = 1186
$ws.Cells.Item(42, 7).Value = 1302

$ws.Cells.Item(43, 6).Value = 1337
$ws.Cells.Item(43, 7).Value = 1352

$ws.Cells.Item(44, 1).Value = 26
$ws.Cells.Item(44, 3).Value = 1444
$ws.Cells.Item(44, 4).Value = 809
$ws.Cells.Item(44, 5).Value = 1296
$ws.Cells.Item(44, 6).Value = 847
$ws.Cells.Item(44, 7).Value = 913

$ws.Cells.Item(45, 6).Value = 972
$ws.Cells.Item(45, 7).Value = 1052

$ws.Cells.Item(46, 6).Value = 1078
$ws.Cells.Item(46, 7).Value = 1081

$ws.Cells.Item(47, 6).Value = 1131
$ws.Cells.Item(47, 7).Value = 1214

$ws.Cells.Item(48, 6).Value = 1242
$ws.Cells.Item(48, 7).Value = 1254

$ws.Cells.Item(49, 1).Value = 28
$ws.Cells.Item(49, 3).Value = 1442
$ws.Cells.Item(49, 4).Value = 862
$ws.Cells.Item(49, 5).Value = 1173
$ws.Cells.Item(49, 6).Value = 889
$ws.Cells.Item(49, 7).Value = 935

$ws.Cells.Item(50, 6).Value = 977
$ws.Cells.Item(50, 7).Value = 984

$ws.Cells.Item(51, 6).Value = 1001
$ws.Cells.Item(51, 7).Value = 1005

$ws.Cells.Item(52, 6).Value = 1018
$ws.Cells.Item(52, 7).Value = 1094

$ws.Cells.Item(53, 6).Value = 1112
$ws.Cells.Item(53, 7).Value = 1133

$ws.Cells.Item(54, 1).Value = 30
$ws.Cells.Item(54, 3).Value = 1461
$ws.Cells.Item(54, 4).Value = 170
$ws.Cells.Item(54, 5).Value = 236
$ws.Cells.Item(54, 6).Value = 191
$ws.Cells.Item(54, 7).Value = 207

$ws.Cells.Item(55, 4).Value = 448
$ws.Cells.Item(55, 5).Value = 537
$ws.Cells.Item(55, 6).Value = 488
$ws.Cells.Item(55, 7).Value = 493

$ws.Cells.Item(56, 4).Value = 1013
$ws.Cells.Item(56, 5).Value = 1300
$ws.Cells.Item(56, 6).Value = 1054
$ws.Cells.Item(56, 7).Value = 1063

$ws.Cells.Item(57, 6).Value = 1104
$ws.Cells.Item(57, 7).Value = 1202

$ws.Cells.Item(58, 6).Value = 1205
$ws.Cells.Item(58, 7).Value = 1258

$ws.Cells.Item(59, 1).Value = 32
$ws.Cells.Item(59, 3).Value = 1460
$ws.Cells.Item(59, 4).Value = 438
$ws.Cells.Item(59, 5).Value = 580
$ws.Cells.Item(59, 6).Value = 473
$ws.Cells.Item(59, 7).Value = 551

$ws.Cells.Item(60, 1).Value = 34
$ws.Cells.Item(60, 3).Value = 1483
$ws.Cells.Item(60, 4).Value = 188
$ws.Cells.Item(60, 5).Value = 270
$ws.Cells.Item(60, 6).Value = 211
$ws.Cells.Item(60, 7).Value = 227

$ws.Cells.Item(61, 4).Value = 403
$ws.Cells.Item(61, 5).Value = 519
$ws.Cells.Item(61, 6).Value = 431
$ws.Cells.Item(61, 7).Value = 478

$ws.Cells.Item(62, 4).Value = 999
$ws.Cells.Item(62, 5).Value = 1186
$ws.Cells.Item(62, 6).Value = 1021
$ws.Cells.Item(62, 7).Value = 1029

$ws.Cells.Item(63, 6).Value = 1076
$ws.Cells.Item(63, 7).Value = 1092

$ws.Cells.Item(64, 6).Value = 1113
$ws.Cells.Item(64, 7).Value = 1125

$ws.Cells.Item(65, 6).Value = 1140
$ws.Cells.Item(65, 7).Value = 1140

$ws.Cells.Item(66, 4).Value = 1320
$ws.Cells.Item(66, 5).Value = 1483
$ws.Cells.Item(66, 6).Value = 1357
$ws.Cells.Item(66, 7).Value = 1364

$ws.Cells.Item(67, 6).Value = 1376
$ws.Cells.Item(67, 7).Value = 1395

$ws.Cells.Item(68, 6).Value = 1448
$ws.Cells.Item(68, 7).Value = 1454

$ws.Cells.Item(69, 1).Value = 36
$ws.Cells.Item(69, 3).Value = 1482
$ws.Cells.Item(69, 4).Value = 288
$ws.Cells.Item(69, 5).Value = 387
$ws.Cells.Item(69, 6).Value = 327
$ws.Cells.Item(69, 7).Value = 344

$ws.Cells.Item(70, 4).Value = 872
$ws.Cells.Item(70, 5).Value = 984
$ws.Cells.Item(70, 6).Value = 905
$ws.Cells.Item(70, 7).Value = 918

$ws.Cells.Item(71, 6).Value = 927
$ws.Cells.Item(71, 7).Value = 944

$ws.Cells.Item(72, 1).Value = 38
$ws.Cells.Item(72, 3).Value = 1492
$ws.Cells.Item(72, 4).Value = 555
$ws.Cells.Item(72, 5).Value = 598
$ws.Cells.Item(72, 6).Value = 650
$ws.Cells.Item(72, 7).Value = 691

$ws.Cells.Item(73, 1).Value = 40
$ws.Cells.Item(73, 3).Value = 1479
$ws.Cells.Item(73, 4).Value = 977
$ws.Cells.Item(73, 5).Value = 1433
$ws.Cells.Item(73, 6).Value = 1007
$ws.Cells.Item(73, 7).Value = 1010

$ws.Cells.Item(74, 6).Value = 1027
$ws.Cells.Item(74, 7).Value = 1036

$ws.Cells.Item(75, 6).Value = 1233
$ws.Cells.Item(75, 7).Value = 1286

$ws.Cells.Item(76, 6).Value = 1295
$ws.Cells.Item(76, 7).Value = 1299

$ws.Cells.Item(77, 6).Value = 1308
$ws.Cells.Item(77, 7).Value = 1319

$ws.Cells.Item(78, 6).Value = 1330
$ws.Cells.Item(78, 7).Value = 1347

$ws.Cells.Item(79, 6).Value = 1357
$ws.Cells.Item(79, 7).Value = 1398

$ws.Cells.Item(80, 1).Value = 42
$ws.Cells.Item(80, 3).Value = 1492
$ws.Cells.Item(80, 4).Value = 906
$ws.Cells.Item(80, 5).Value = 1343
$ws.Cells.Item(80, 6).Value = 930
$ws.Cells.Item(80, 7).Value = 952

$ws.Cells.Item(81, 6).Value = 960
$ws.Cells.Item(81, 7).Value = 998

$ws.Cells.Item(82, 6).Value = 1034
$ws.Cells.Item(82, 7).Value = 1129

$ws.Cells.Item(83, 6).Value = 1148
$ws.Cells.Item(83, 7).Value = 1186

$ws.Cells.Item(84, 6).Value = 1209
$ws.Cells.Item(84, 7).Value = 1248

$ws.Cells.Item(85, 6).Value = 1293
$ws.Cells.Item(85, 7).Value = 1300

$ws.Cells.Item(86, 1).Value = 44
$ws.Cells.Item(86, 3).Value = 1481
$ws.Cells.Item(86, 4).Value = 285
$ws.Cells.Item(86, 5).Value = 380
$ws.Cells.Item(86, 6).Value = 318
$ws.Cells.Item(86, 7).Value = 353

$ws.Cells.Item(87, 4).Value = 948
$ws.Cells.Item(87, 5).Value = 1230
$ws.Cells.Item(87, 6).Value = 985
$ws.Cells.Item(87, 7).Value = 994

$ws.Cells.Item(88, 6).Value = 1027
$ws.Cells.Item(88, 7).Value = 1035

$ws.Cells.Item(89, 6).Value = 1045
$ws.Cells.Item(89, 7).Value = 1048

$ws.Cells.Item(90, 1).Value = 45
$ws.Cells.Item(90, 3).Value = 1497
$ws.Cells.Item(90, 4).Value = 1068
$ws.Cells.Item(90, 5).Value = 1347
$ws.Cells.Item(90, 6).Value = 1094
$ws.Cells.Item(90, 7).Value = 1100

$ws.Cells.Item(91, 6).Value = 1166
$ws.Cells.Item(91, 7).Value = 1199

$ws.Cells.Item(92, 6).Value = 1204
$ws.Cells.Item(92, 7).Value = 1286

# Move the active selection to the last entry added (orbit 45 / row 92),
# matching where data entry left off.
$ws.Range("C92").Select()
